$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): strip the bold/border/center-top-aligned style that
# was previously applied to A1:AR1, and clear the stray "Unnamed: 0" label
# in A1 so it matches the rest of the (now-blank) label cells. ---
$headerRange = $ws.Range("A1:AR1")
$headerRange.Style = "Normal"
$ws.Range("A1").Value = ""

# --- Corrected "summary" (AH) column totals for rows 4, 5, 6, 8 ---
$ws.Range("AH4").Value = 59
$ws.Range("AH5").Value = 163
$ws.Range("AH6").Value = 85757.2
$ws.Range("AH8").Value = 526.12

# --- Row 7 ("Dwell time (%)") recomputed against the corrected totals ---
$ws.Range("B7").Value = 0.99
$ws.Range("C7").Value = 0.5600000000000001
$ws.Range("D7").Value = 0.08
$ws.Range("E7").Value = 15.32
$ws.Range("F7").Value = 4.97
$ws.Range("G7").Value = 2
$ws.Range("H7").Value = 1.09
$ws.Range("I7").Value = 0.13
$ws.Range("J7").Value = 0.95
$ws.Range("L7").Value = 0.71
$ws.Range("M7").Value = 1.12
$ws.Range("N7").Value = 0.74
$ws.Range("O7").Value = 0.46
$ws.Range("Q7").Value = 0.08
$ws.Range("R7").Value = 0.66
$ws.Range("S7").Value = 0.44
$ws.Range("U7").Value = 0.26
$ws.Range("W7").Value = 0.74
$ws.Range("X7").Value = 0.06
$ws.Range("Z7").Value = 0.24
$ws.Range("AA7").Value = 0.06
$ws.Range("AC7").Value = 0.15
$ws.Range("AF7").Value = 3.04
$ws.Range("AH7").Value = 35.75
$ws.Range("AI7").Value = 0.15
$ws.Range("AJ7").Value = 0.08
$ws.Range("AK7").Value = 0.18
$ws.Range("AL7").Value = 0.34
$ws.Range("AN7").Value = 0.07000000000000001
$ws.Range("AO7").Value = 0.74
$ws.Range("AP7").Value = 0.06

# --- Drop the two trailing blank rows (11 and 12); Excel auto-shrinks the
# sheet dimension from A1:AR12 to A1:AR10 as a result. ---
$ws.Rows("11:12").Delete()
